$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the "ID GENERATO" values previously filled in column S for rows 3-5
$ws.Range("S3").Value = ""
$ws.Range("S4").Value = ""
$ws.Range("S5").Value = ""
